# Update existing rows: replace the `""` / `empty list` / `[] ` values
# in the "Expected Output" / "Actual Output" columns with "enter string 2-9"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "enter string 2-9"
$ws.Range("E4").Value = "enter string 2-9"

$ws.Range("D5").Value = "enter string 2-9"
$ws.Range("E5").Value = "enter string 2-9"

$ws.Range("D6").Value = "enter string 2-9"
$ws.Range("E6").Value = "enter string 2-9"

# Add a new test case row (row 7) for "with expression"
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "with expression"
$ws.Range("C7").Value = "2+3"
$ws.Range("D7").Value = "enter string 2-9"
$ws.Range("E7").Value = "enter string 2-9"
$ws.Range("F7").Value = "FAIL"

# Move the active selection down to follow the newly added row
$ws.Range("F8").Select() | Out-Null
